$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.314.29'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.30%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.853.90'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.45%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.93%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.32'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.93%  '

$ws.Range('E6').Value = '  -0.82%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4613'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.79%  '

$ws.Range('E8').Value = '  +0.26%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07312'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.67%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8849'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.39%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.94'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.07%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07785'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.62%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.849.57'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +6.10%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.384'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.71%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.549'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.48%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.71'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.19%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.003'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.81%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009029'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.36%  '

$ws.Range('E19').Value = '  -0.80%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.29%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.326.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.42%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.135'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.10%  '

$ws.Range('E23').Value = '  -0.48%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.073.95'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.05%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.928'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.96%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.51%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.51%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.073'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.11%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.103'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.09%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.03'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.86%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08853'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.13%  '

$ws.Range('B32').Value = 'HuobiToken'
$ws.Range('C32').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.120'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.09%  '

$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7779'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +7.52%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.173'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.54%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.504'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.36%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.674'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.43%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01961'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.77%  '

$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.078'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.35%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05234'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.17%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.966'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.52%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.982'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.57%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5145'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.62%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1633'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.38%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.420'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.94%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4816'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.00%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.34'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.70%  '

$ws.Range('E47').Value = '  -0.92%  '

$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.649'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.51%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '102.37'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.07%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06217'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.04%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '65.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.33%  '
